$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Expand the "RS" / "BM" abbreviations that appear together in the
# second paragraph's trailing run ("...multiple devices on the RS side. In
# this case when someone speaks at BM side..."). Using the combined, unique
# phrase (rather than "RS side" / "BM side" separately) avoids touching the
# earlier, unrelated "BM side" occurrence that appears in the very first run
# of that same paragraph ("...speaks at the BM side, the application...").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "RS side. In this case when someone speaks at BM side",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Receiver Side side. In this case when someone speaks at Broadcaster Microphone side",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Step 2: Append a brand-new paragraph (Times New Roman / 28pt, same as the
# surrounding body text) describing the "Second Approach". A short marker is
# appended after the real text so the trailing "_GoBack" bookmark can later
# be re-anchored at a position that is not the literal last character of the
# document/story (collapsed bookmarks placed at the very last offset land in
# the wrong spot in this host), and then the marker text is removed again.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastRange.Collapse(0)
$lastRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Collapse(0)

$secondApproachText = "Second Approach: In this stage a device would take audio with noises from nearby sources and produce better sound quality as output in the other device. Initially it would be done with two devices. Usually modern cellphones have numerous microphones. So, we can assume that a device would take sound input from its microphones and broadcast the sound with better quality audio to a nearby device which would be automatically connected when bought into the range of the broadcasting device. After this approach is successfully done, we can try with numerous android devices which would take audio input simultaneously and broadcast to another android phone. Here the number of cellphones in the receiving side may vary with project progress."
$endMarker = "@@END_MARK@@"
$newRange.InsertAfter($secondApproachText + $endMarker)

# ---------------------------------------------------------------------------
# Step 3: Re-anchor the hidden "_GoBack" bookmark so it sits, collapsed, at
# the very end of the document (i.e. the end of the paragraph just added) -
# matching Word's normal behaviour of moving _GoBack to the last edit point.
# ---------------------------------------------------------------------------
$markerRange = $d.Content
$markerRange.Find.Execute($endMarker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markerRange.Collapse(1) | Out-Null
$d.Bookmarks.Add("_GoBack", $markerRange) | Out-Null

# Remove the temporary end marker text again.
$d.Content.Find.Execute($endMarker, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

Write-Host "Edit complete"
